# show_command.xlsx edit: "change port and add xr"
#
# Semantics (derived from the canonical OOXML diff):
#   - A new device-type row for "cisco_xr" is inserted at row 6. Its command
#     set (columns B:DD) is a copy of the "cisco_xe" row (row 2), since IOS XR
#     shares nearly all commands with IOS XE in this sheet.
#   - The row that used to live at row 6 ("riverbed", columns A:AQ) moves down
#     to row 7, with two obsolete commands dropped from it:
#       "show hardware error-log" and "show in-path neighbour"
#     (the remaining cells shift left to stay contiguous).
#   - The active-cell selection moves from E11 to C11 ("change port").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Capture the current ("riverbed") row 6 contents before it gets overwritten.
$oldRow6 = $ws.Range("A6:AQ6").Value()

# 2) Commands that no longer apply to the riverbed row once it shifts to row 7.
$dropped = @("show hardware error-log", "show in-path neighbour")

$kept = @()
for ($i = 1; $i -le 43; $i++) {
    $val = $oldRow6[1, $i]
    if (-not ($dropped -contains $val)) {
        $kept += $val
    }
}

# 3) Write the filtered riverbed data into row 7.
for ($i = 0; $i -lt $kept.Count; $i++) {
    $ws.Cells.Item(7, $i + 1).Value = $kept[$i]
}

# 4) Turn row 6 into the new "cisco_xr" row: label in column A, and the
#    cisco_xe command set (row 2) copied across columns B:DD.
$ws.Range("A6").Value = "cisco_xr"
$ws.Range("B6:DD6").Value = $ws.Range("B2:DD2").Value()

# 5) Move the active selection from E11 to C11.
$ws.Range("C11").Select()
